# Generate Report for Handoff
# Updates the localization-status workbook to reflect a new handoff:
#   - Status changes from "In Translation" to "Ready for handoff"
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps advance
#   - The "Status"/date columns on all three sheets are widened to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-17 06:34:37"

$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-17 06:34:32"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-17 06:34:37"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
